$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2..6), matching columns D, L, M, N, O, P, S
$data = @{
    2 = @{ D = 44204; L = "Primera";  M = 110; N = 7000;  O = 7500;  P = 7318;  S = 1045 }
    3 = @{ D = 44189; L = "Especial"; M = 20;  N = 15000; O = 15000; P = 15000; S = 2143 }
    4 = @{ D = 44187; L = "Primera";  M = 30;  N = 13000; O = 13000; P = 13000; S = 1857 }
    5 = @{ D = 44187; L = "Especial"; M = 45;  N = 14000; O = 14000; P = 14000; S = 2000 }
    6 = @{ D = 44187; L = "Primera";  M = 50;  N = 12000; O = 12000; P = 12000; S = 1714 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
